$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying numeric-looking values in columns C (runs), D (balls) and
# E (fours) are stored as text, not numbers (see ignoredErrors/numberStoredAsText
# on A1:F4 in the original workbook). Force the target ranges to a text
# number format before assigning so Excel keeps them as text instead of
# auto-converting them to numeric values.
$ws.Range("C2:E4").NumberFormat = "@"

# Row 2
$ws.Range("C2").Value = "14"
$ws.Range("D2").Value = "9"
$ws.Range("E2").Value = "2"

# Row 3
$ws.Range("C3").Value = "7"
$ws.Range("D3").Value = "8"
$ws.Range("E3").Value = "1"

# Row 4
$ws.Range("C4").Value = "12"
$ws.Range("D4").Value = "12"
# E4 stays "2" (unchanged)
